$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Replace the StatQuery (column C) text on rows 2-4 with the new Cypher
#    query (same text is shared by all three rows, as in the original file).
# ---------------------------------------------------------------------------
$newQuery = @"
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['Boxer']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS ``Case Files``,
    count(distinct sf) AS ``Study Files``
"@

$ws.Range("C2").Value2 = $newQuery
$ws.Range("C3").Value2 = $newQuery
$ws.Range("C4").Value2 = $newQuery

# ---------------------------------------------------------------------------
# 2. Row heights shrink now that the new query text is shorter (previously
#    they were all pinned at Excel's max of 409.6pt).
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 201.6
$ws.Rows.Item(3).RowHeight = 230.4
$ws.Rows.Item(4).RowHeight = 244.8

# ---------------------------------------------------------------------------
# 3. Column widths were re-fit to the new content. The host's column-width
#    model quantizes to 1/6-character pixel steps (same granularity Excel
#    itself uses for ColumnWidth), so the inputs below are chosen as the
#    closest values that round-trip to the target widths on save.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 10.0
$ws.Columns.Item(2).ColumnWidth = 78.66666666666667
$ws.Columns.Item(3).ColumnWidth = 58.0
$ws.Columns.Item(4).ColumnWidth = 41.333333333333336
$ws.Columns.Item(5).ColumnWidth = 40.166666666666664

# ---------------------------------------------------------------------------
# 4. Selection / scroll position moved from B3 to B4, with the view
#    scrolled back up so row 1 is visible again (topLeftCell B1 vs B3).
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 2
[void]$ws.Range("B4").Select()
